$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1755.25
$ws.Range("I20").Value = 673.6667
$ws.Range("K20").Value = 673.6667
$ws.Range("M20").Value = -443.6667

$ws.Range("H26").Value = 1000
$ws.Range("I26").Value = 1000
$ws.Range("K26").Value = 1000
$ws.Range("M26").Value = -656

$ws.Range("H35").Value = 1755.25
$ws.Range("I35").Value = 673.6667
$ws.Range("K35").Value = 673.6667
$ws.Range("M35").Value = -294.6667

$ws.Range("H55").Value = 744.9167
$ws.Range("J55").Value = 1179.6
$ws.Range("L55").Value = 1179.6
$ws.Range("N55").Value = -1607.6

$ws.Range("H129").Value = 2919.8333
$ws.Range("J129").Value = 5000
$ws.Range("L129").Value = 15000
$ws.Range("N129").Value = -25000

$ws.Range("H135").Value = 2101.1667
$ws.Range("I135").Value = 1672.1111
$ws.Range("K135").Value = 15048.9999
$ws.Range("M135").Value = -12513.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 3669
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 3003.5
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 3003.5
$ws.Range("M16").Value = -4713
$ws.Range("N16").Value = -3577.5

$ws.Range("H32").Value = 4863.8438
$ws.Range("I32").Value = 3430.7778
$ws.Range("K32").Value = 3430.7778
$ws.Range("M32").Value = -3143.7778

$ws.Range("H45").Value = 2999.6
$ws.Range("I45").Value = 1999
$ws.Range("J45").Value = 3666.6667
$ws.Range("K45").Value = 1999
$ws.Range("L45").Value = 3666.6667
$ws.Range("M45").Value = -1622
$ws.Range("N45").Value = -4420.6667

$ws.Range("H74").Value = 3258.9412
$ws.Range("I74").Value = 3027.1333
$ws.Range("K74").Value = 3027.1333
$ws.Range("M74").Value = -2153.1333

$ws.Range("H77").Value = 3258.9412
$ws.Range("I77").Value = 3027.1333
$ws.Range("K77").Value = 15135.6665
$ws.Range("M77").Value = -10767.6665

$ws.Range("H102").Value = 10103186
$ws.Range("I102").Value = 15874550
$ws.Range("K102").Value = 15874550
$ws.Range("M102").Value = -15872928

$ws.Range("H110").Value = 3833576.8
$ws.Range("I110").Value = 6174329
$ws.Range("K110").Value = 6174329
$ws.Range("M110").Value = -6172284

$ws.Range("H132").Value = 3459.4167
$ws.Range("I132").Value = 3390.611
$ws.Range("K132").Value = 10171.833
$ws.Range("M132").Value = -7641.832999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4914.625
$ws.Range("J20").Value = 6581.5
$ws.Range("L20").Value = 6581.5
$ws.Range("N20").Value = -7075.5

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H75").Value = 1000
$ws.Range("I75").Value = 1000
$ws.Range("K75").Value = 1000
$ws.Range("M75").Value = -64

$ws.Range("H78").Value = 1000
$ws.Range("I78").Value = 1000
$ws.Range("K78").Value = 3000
$ws.Range("M78").Value = 1680

$ws.Range("H99").Value = 2509.1428
$ws.Range("I99").Value = 2509.1428
$ws.Range("K99").Value = 2509.1428
$ws.Range("M99").Value = -1011.1428

$ws.Range("H105").Value = 3090096.8
$ws.Range("I105").Value = 5211404
$ws.Range("J105").Value = 4558.636
$ws.Range("K105").Value = 5211404
$ws.Range("L105").Value = 4558.636
$ws.Range("M105").Value = -5209657
$ws.Range("N105").Value = -8052.636

$ws.Range("H134").Value = 1305.6
$ws.Range("I134").Value = 1367.5555
$ws.Range("J134").Value = 748
$ws.Range("K134").Value = 4102.666499999999
$ws.Range("L134").Value = 2244
$ws.Range("M134").Value = -1567.666499999999
$ws.Range("N134").Value = -7314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 3901.25
$ws.Range("I15").Value = 7357.5
$ws.Range("K15").Value = 7357.5
$ws.Range("M15").Value = -7187.5

$ws.Range("H21").Value = 3000
$ws.Range("J21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3470

$ws.Range("H25").Value = 4950
$ws.Range("I25").Value = 2425
$ws.Range("K25").Value = 2425
$ws.Range("M25").Value = -2251

$ws.Range("H26").Value = 15000
$ws.Range("I26").Value = 15000
$ws.Range("K26").Value = 15000
$ws.Range("M26").Value = -14713

$ws.Range("H31").Value = 3359.8064
$ws.Range("I31").Value = 2223.2856
$ws.Range("J31").Value = 5746.5
$ws.Range("K31").Value = 2223.2856
$ws.Range("L31").Value = 5746.5
$ws.Range("M31").Value = -1928.2856
$ws.Range("N31").Value = -6336.5

$ws.Range("H34").Value = 3359.8064
$ws.Range("I34").Value = 2223.2856
$ws.Range("J34").Value = 5746.5
$ws.Range("K34").Value = 2223.2856
$ws.Range("L34").Value = 5746.5
$ws.Range("M34").Value = -2021.2856
$ws.Range("N34").Value = -6150.5

$ws.Range("H58").Value = 3833.3333
$ws.Range("I58").Value = 2003.6666
$ws.Range("J58").Value = 4748.1665
$ws.Range("K58").Value = 2003.6666
$ws.Range("L58").Value = 4748.1665
$ws.Range("M58").Value = -1800.6666
$ws.Range("N58").Value = -5154.1665

$ws.Range("H86").Value = 6547.5713
$ws.Range("I86").Value = 6547.5713
$ws.Range("K86").Value = 6547.5713
$ws.Range("M86").Value = -5424.5713

$ws.Range("H89").Value = 6547.5713
$ws.Range("I89").Value = 6547.5713
$ws.Range("K89").Value = 32737.8565
$ws.Range("M89").Value = -27121.8565

$ws.Range("H99").Value = 12832.926
$ws.Range("I99").Value = 9348.909
$ws.Range("K99").Value = 9348.909
$ws.Range("M99").Value = -7850.909

$ws.Range("H105").Value = 2999.818
$ws.Range("I105").Value = 1832.6666
$ws.Range("K105").Value = 1832.6666
$ws.Range("M105").Value = -85.66660000000002

$ws.Range("H107").Value = 13889819
$ws.Range("I107").Value = 23810050
$ws.Range("K107").Value = 23810050
$ws.Range("M107").Value = -23808130

$ws.Range("H126").Value = 12832.926
$ws.Range("I126").Value = 9348.909
$ws.Range("K126").Value = 28046.727
$ws.Range("M126").Value = -25576.727

$ws.Range("H136").Value = 3833.3333
$ws.Range("I136").Value = 2003.6666
$ws.Range("J136").Value = 4748.1665
$ws.Range("K136").Value = 6010.9998
$ws.Range("L136").Value = 14244.4995
$ws.Range("M136").Value = -3460.9998
$ws.Range("N136").Value = -19344.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 194
$ws.Range("I8").Value = 194
$ws.Range("K8").Value = 582
$ws.Range("M8").Value = -443

$ws.Range("H23").Value = 999999
$ws.Range("J23").Value = 999999
$ws.Range("L23").Value = 2999997
$ws.Range("N23").Value = -3000467

$ws.Range("H131").Value = 1377.2
$ws.Range("I131").Value = 608.5714
$ws.Range("J131").Value = 1478.717
$ws.Range("K131").Value = 1825.7142
$ws.Range("L131").Value = 4436.151
$ws.Range("M131").Value = 3214.2858
$ws.Range("N131").Value = -14516.151

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 32499.5
$ws.Range("J51").Value = 34999
$ws.Range("L51").Value = 34999
$ws.Range("N51").Value = -36017

$ws.Range("H132").Value = 2995
$ws.Range("I132").Value = 2995
$ws.Range("K132").Value = 8985
$ws.Range("M132").Value = -6455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4008
$ws.Range("I4").Value = 4009
$ws.Range("K4").Value = 4009
$ws.Range("M4").Value = -3896

$ws.Range("H28").Value = 4008
$ws.Range("I28").Value = 4009
$ws.Range("K28").Value = 4009
$ws.Range("M28").Value = -3777

$ws.Range("H34").Value = 6062
$ws.Range("I34").Value = 9000
$ws.Range("J34").Value = 3124
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 3124
$ws.Range("M34").Value = -8828
$ws.Range("N34").Value = -3468

$ws.Range("H37").Value = 4008
$ws.Range("I37").Value = 4009
$ws.Range("K37").Value = 4009
$ws.Range("M37").Value = -3902

$ws.Range("H41").Value = 33999
$ws.Range("I41").Value = 33999
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 33999
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -33561
$ws.Range("N41").ClearContents()

$ws.Range("H46").Value = 2247.5557
$ws.Range("I46").Value = 995.2
$ws.Range("J46").Value = 2729.2307
$ws.Range("K46").Value = 995.2
$ws.Range("L46").Value = 2729.2307
$ws.Range("M46").Value = -807.2
$ws.Range("N46").Value = -3105.2307

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H55").Value = 523.2632
$ws.Range("I55").Value = 453.7143
$ws.Range("J55").Value = 718
$ws.Range("K55").Value = 453.7143
$ws.Range("L55").Value = 718
$ws.Range("M55").Value = -280.7143
$ws.Range("N55").Value = -1064

$ws.Range("H61").Value = 55558056
$ws.Range("I61").Value = 55558056
$ws.Range("K61").Value = 55558056
$ws.Range("M61").Value = -55557854

$ws.Range("H100").Value = 5600.8
$ws.Range("I100").Value = 5600.8
$ws.Range("K100").Value = 5600.8
$ws.Range("M100").Value = -5059.8

$ws.Range("H113").Value = 55558056
$ws.Range("I113").Value = 55558056
$ws.Range("K113").Value = 55558056
$ws.Range("M113").Value = -55555886

$ws.Range("H136").Value = 4906.5884
$ws.Range("I136").Value = 4425
$ws.Range("K136").Value = 13275
$ws.Range("M136").Value = -10725
